# Word COM-interop script applying the "Biosphere Engineering" -> "Exploring Music"
# rewrite described by the target diff.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "WARN: replace failed for: $old"
    }
}

function Find-End([string]$needle) {
    # Returns the document character offset just after $needle (a fresh, non-aliased
    # Range object is used for the actual match so later collapses/inserts don't
    # get confused by Range auto-reseating).
    $rng = $d.Content
    $ok = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Output "WARN: find failed for: $needle"
        return -1
    }
    return $rng.End
}

# ---------------------------------------------------------------------------
# Title / byline / contact block
# ---------------------------------------------------------------------------
Replace-Text "Biosphere Engineering: Shaping Life's Course" "Exploring Music: Harmonizing Sounds and Emotions"
Replace-Text "Sophia Oliver" "Isabella Clark"
Replace-Text "sophia" "isabella"
Replace-Text "oliver@biosphereengineering" "clark@education"
Replace-Text "org" "com"

# ---------------------------------------------------------------------------
# Body paragraph 1 (intro)
# ---------------------------------------------------------------------------
Replace-Text "In the vast and intricate tapestry of life, humanity stands at a pivotal juncture, poised to wield immense power over the very fabric of existence" "Music, a universal language that transcends borders and cultures, is a captivating force in our lives"

Replace-Text " Biosphere engineering, an emerging field of transformative potential, beckons us to question the boundaries between the natural and the artificial, inviting us to co-create a sustainable and harmonious coexistence with our planet" " It permeates our emotions, influences our moods, and holds the power to transport us to different realms"

Replace-Text " As we grapple with the repercussions of our actions on the Earth's intricate ecosystems, this field offers a beacon of hope, promising a comprehensive, ethical, and holistic approach to nurturing life's flourishing" " The synergy between sounds, rhythms, and melodies weaves a rich tapestry, inviting us to delve into its intricacies and discover its mesmerizing allure"

$p = Find-End "its mesmerizing allure"
if ($p -ge 0) {
    $d.Range($p, $p).InsertAfter(". Join us on a harmonious journey as we explore the multifaceted world of music, unveiling its ability to evoke emotions, narrate stories, and connect humanity")
}

Replace-Text "Life on Earth has gracefully unfolded over eons, weaving a symphony of interconnectedness" "In the world of music, the connection between sounds and emotions is profound"

Replace-Text " Through biosphere engineering, we possess the power to influence this dance of life, crafting a narrative of symbiosis and resilience" " Certain melodies, harmonies, and rhythms can trigger specific emotional responses, ranging from joy and exuberance to sadness and nostalgia"

Replace-Text " By delving into the underlying principles of ecosystems, we gain the ability to mindfully guide species interactions, enhancing biodiversity, ensuring food security, and mitigating the impacts of climate change" " This phenomenon, known as the psychology of music, has been extensively studied, revealing the intricate relationship between auditory stimuli and human emotions"

$p = Find-End "auditory stimuli and human emotions"
if ($p -ge 0) {
    $d.Range($p, $p).InsertAfter(". Music has the uncanny ability to tap into our deepest feelings and resonate with our experiences, transporting us to a realm where emotions flow freely")
}

Replace-Text "With great power comes immense responsibility" "Music is an art form that captures and reflects the human condition"

Replace-Text " Biosphere engineering demands an unwavering commitment to humility, wisdom, and sustainability" " It narrates stories of love, loss, triumph, and despair, weaving tales that mirror the ebb and flow of life"

Replace-Text " It necessitates an ethical framework that honors the intrinsic value of all living beings and safeguards the delicate balance of natural systems" " Through lyrics and melodies, music conveys messages that transcend words, allowing us to connect with experiences beyond our own"

Replace-Text " As we embark on this journey, interdisciplinary collaboration and unwavering curiosity shall be our guiding stars, illuminating our path towards flourishing ecosystems and a vibrant future for life on Earth" " It captures the zeitgeist of an era, encapsulating the hopes, dreams, and fears of a generation"

$p = Find-End "fears of a generation"
if ($p -ge 0) {
    $target = $d.Range($p, $p)
    $target.InsertAfter(". Music becomes a soundtrack to our lives, accompanying us through milestones and marking significant moments with its poignant melodies.")
    $target2 = $d.Range($target.End, $target.End)
    $target2.InsertAfter([char]11)
    $target3 = $d.Range($target2.End, $target2.End)
    # Note: real Word would stamp a <w:lastRenderedPageBreak/> rendering-cache marker
    # right before "and mourning" once it repaginates (that's where the page happens
    # to break after this edit). That marker is written by Word's layout engine only
    # -- it has no corresponding writable property/method on the Find/Range object
    # model, so headless COM automation cannot author it directly.
    $target3.InsertAfter([char]11 + "Music transcends cultural and geographical boundaries, uniting humanity in a shared experience. It serves as a bridge between people, fostering understanding and appreciation for diverse cultures. Through shared melodies and rhythms, music creates a sense of community, bringing people together in moments of celebration, worship, and mourning. It celebrates our common humanity, reminding us that despite our differences, we are all connected by the universal language of music")
}

# ---------------------------------------------------------------------------
# Summary paragraph
# ---------------------------------------------------------------------------
Replace-Text "Biosphere engineering emerges as a transformative field, empowering humanity to shape the course of life on Earth" "In this exploration of music, we have delved into its ability to evoke emotions, narrate stories, and connect humanity"

Replace-Text " With the potential to influence species interactions, enhance biodiversity, ensure food security, and mitigate climate change, this discipline offers a comprehensive approach to nurturing ecological resilience" " Music is a powerful medium that communicates emotions that words cannot express, creating a tapestry of sound that resonates with our souls"

Replace-Text " Guided by humility, wisdom, and sustainability, biosphere engineering calls for an unwavering commitment to ethical decision-making and interdisciplinary collaboration" " It serves as a mirror to society, reflecting our joys, sorrows, and collective experiences"

Replace-Text " This field invites us to co-create a sustainable and harmonious coexistence with our planet, shaping life's journey towards a flourishing future" " Moreover, music transcends cultural barriers, fostering unity and understanding among people from all walks of life"

$p = Find-End "from all walks of life"
if ($p -ge 0) {
    $d.Range($p, $p).InsertAfter(". As a universal language, it brings humanity together, creating a harmonious symphony that celebrates our shared existence")
}

# ---------------------------------------------------------------------------
# Trailing empty paragraph
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count
$last = $d.Paragraphs($n).Range
$last.InsertParagraphAfter()
